# Updated cryptos list (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.792.30"
$ws.Range("E2").Value = "  -1.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.352.87"
$ws.Range("E3").Value = "  -2.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.27%  "

# Row 5 - BNB
$ws.Range("D5").Value = "320.91"
$ws.Range("E5").Value = "  -2.09%  "

# Row 6 - Solana
$ws.Range("D6").Value = "105.37"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  -2.79%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  -7.43%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "41.26"
$ws.Range("E10").Value = "  -2.51%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0924"
$ws.Range("E11").Value = "  -2.28%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "8.44"
$ws.Range("E12").Value = "  -2.37%  "

# Row 13 - Polygon
$ws.Range("E13").Value = "  -2.33%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.07%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "15.99"
$ws.Range("E15").Value = "  -7.64%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.706.33"
$ws.Range("E16").Value = "  -2.50%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.310.28"
$ws.Range("E17").Value = "  -4.59%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.763.85"
$ws.Range("E18").Value = "  -1.84%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +2.69%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -3.37%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "77.17"
$ws.Range("E21").Value = "  +1.30%  "

# Row 22 - PancakeSwap
$ws.Range("D22").Value = "3.62"
$ws.Range("E22").Value = "  +2.98%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "261.55"
$ws.Range("E23").Value = "  -3.93%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "2.33"
$ws.Range("E24").Value = "  -5.32%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "9.60"
$ws.Range("E25").Value = "  -0.40%  "

# Row 26 - Dai: unchanged

# Row 27 - Cosmos
$ws.Range("D27").Value = "11.39"
$ws.Range("E27").Value = "  -5.26%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "23.22"
$ws.Range("E28").Value = "  +0.69%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -0.01%  "

# Row 30 - Monero
$ws.Range("D30").Value = "174.86"
$ws.Range("E30").Value = "  -1.72%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "36.34"
$ws.Range("E31").Value = "  -4.08%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "6.16"
$ws.Range("E32").Value = "  +3.27%  "

# Rows 33/34 swapped order: Hedera now ranked above WEMIXToken
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0890"
$ws.Range("E33").Value = "  -5.22%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.99"
$ws.Range("E34").Value = "  -7.19%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  -3.09%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  +8.31%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "4.63"
$ws.Range("E37").Value = "  -5.17%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.0361"
$ws.Range("E38").Value = "  -2.96%  "

# Row 39 - NEARProtocol
$ws.Range("D39").Value = "3.81"
$ws.Range("E39").Value = "  -6.83%  "

# Row 40 - LidoDAOToken
$ws.Range("D40").Value = "2.69"
$ws.Range("E40").Value = "  -6.51%  "

# Row 41 - MultiversX
$ws.Range("E41").Value = "  +2.35%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  -9.48%  "

# Row 43 - Algorand
$ws.Range("D43").Value = "0.232"
$ws.Range("E43").Value = "  -1.71%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.26%  "

# Row 45 - Aave
$ws.Range("D45").Value = "115.44"
$ws.Range("E45").Value = "  -9.64%  "

# Row 46 - BitcoinSV
$ws.Range("D46").Value = "88.79"
$ws.Range("E46").Value = "  +1.23%  "

# Row 47 - Celestia
$ws.Range("D47").Value = "11.91"
$ws.Range("E47").Value = "  -6.87%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "5.50"
$ws.Range("E48").Value = "  -3.63%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  -6.14%  "

# Row 50 - ordi
$ws.Range("D50").Value = "73.41"
$ws.Range("E50").Value = "  +0.31%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  -4.71%  "
